# Auto-generated edit script applying numeric updates described in the commit diff.
# Each hunk corresponds to a single row in one of the 8 job worksheets (ALC, ARM, BSM,
# CRP, CUL, GSM, LTW, WVR). Values in columns H-N are updated; a couple of cells are
# newly populated (previously blank) and a couple are cleared (previously populated).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 436.57144
$ws.Range("I2").Value = 436.57144
$ws.Range("K2").Value = 436.57144
$ws.Range("M2").Value = -323.57144

$ws.Range("H88").Value = 7034.9165
$ws.Range("I88").Value = 12115.556
$ws.Range("J88").Value = 3986.5334
$ws.Range("K88").Value = 12115.556
$ws.Range("L88").Value = 3986.5334
$ws.Range("M88").Value = -11709.556
$ws.Range("N88").Value = -4798.5334

$ws.Range("H91").Value = 7034.9165
$ws.Range("I91").Value = 12115.556
$ws.Range("J91").Value = 3986.5334
$ws.Range("K91").Value = 12115.556
$ws.Range("L91").Value = 3986.5334
$ws.Range("M91").Value = -10711.556
$ws.Range("N91").Value = -6794.5334

$ws.Range("H107").Value = 35933.332
$ws.Range("I107").Value = 35933.332
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 35933.332
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -34013.332
$ws.Range("N107").ClearContents()

$ws.Range("H111").Value = 1497.5
$ws.Range("I111").Value = 995
$ws.Range("J111").Value = 1665
$ws.Range("K111").Value = 2985
$ws.Range("L111").Value = 4995
$ws.Range("M111").Value = 82
$ws.Range("N111").Value = -11129

$ws.Range("H121").Value = 2119
$ws.Range("J121").Value = 2119
$ws.Range("L121").Value = 6357
$ws.Range("N121").Value = -9851

$ws.Range("H132").Value = 5170.5713
$ws.Range("I132").Value = 5507.706
$ws.Range("K132").Value = 16523.118
$ws.Range("M132").Value = -13993.118

$ws.Range("H138").Value = 9395.366
$ws.Range("I138").Value = 12124.875
$ws.Range("J138").Value = 8733.666999999999
$ws.Range("K138").Value = 36374.625
$ws.Range("L138").Value = 26201.001
$ws.Range("M138").Value = -31234.625
$ws.Range("N138").Value = -36481.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 60071.895
$ws.Range("I2").Value = 12874.777
$ws.Range("K2").Value = 12874.777
$ws.Range("M2").Value = -12761.777

$ws.Range("H32").Value = 1551.2949
$ws.Range("I32").Value = 1578.581
$ws.Range("K32").Value = 1578.581
$ws.Range("M32").Value = -1291.581

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws.Range("H45").Value = 6922.613
$ws.Range("I45").Value = 11137.3125
$ws.Range("J45").Value = 2426.9333
$ws.Range("K45").Value = 11137.3125
$ws.Range("L45").Value = 2426.9333
$ws.Range("M45").Value = -10760.3125
$ws.Range("N45").Value = -3180.9333

$ws.Range("H55").Value = 56682.332
$ws.Range("I55").Value = 30048
$ws.Range("K55").Value = 30048
$ws.Range("M55").Value = -29733

$ws.Range("H74").Value = 49096.566
$ws.Range("I74").Value = 52264.18
$ws.Range("K74").Value = 52264.18
$ws.Range("M74").Value = -51390.18

$ws.Range("H77").Value = 49096.566
$ws.Range("I77").Value = 52264.18
$ws.Range("K77").Value = 261320.9
$ws.Range("M77").Value = -256952.9

$ws.Range("H97").Value = 6670746.5
$ws.Range("I97").Value = 4591.8
$ws.Range("K97").Value = 4591.8
$ws.Range("M97").Value = -4095.8

$ws.Range("H110").Value = 1516.4445
$ws.Range("I110").Value = 716.36365
$ws.Range("K110").Value = 716.36365
$ws.Range("M110").Value = 1328.63635

$ws.Range("H116").Value = 60071.895
$ws.Range("I116").Value = 12874.777
$ws.Range("K116").Value = 12874.777
$ws.Range("M116").Value = -10580.777

$ws.Range("H122").Value = 252333.28
$ws.Range("I122").Value = 1724.3143
$ws.Range("K122").Value = 5172.9429
$ws.Range("M122").Value = -2722.9429

$ws.Range("H132").Value = 2620.5833
$ws.Range("I132").Value = 2238.9375
$ws.Range("K132").Value = 6716.8125
$ws.Range("M132").Value = -4186.8125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 60071.895
$ws.Range("I3").Value = 12874.777
$ws.Range("K3").Value = 12874.777
$ws.Range("M3").Value = -12760.777

$ws.Range("H5").Value = 1673333.4
$ws.Range("J5").Value = 1673333.4
$ws.Range("L5").Value = 1673333.4
$ws.Range("N5").Value = -1673559.4

$ws.Range("H105").Value = 10685.883
$ws.Range("I105").Value = 13611.091
$ws.Range("J105").Value = 5323
$ws.Range("K105").Value = 13611.091
$ws.Range("L105").Value = 5323
$ws.Range("M105").Value = -11864.091
$ws.Range("N105").Value = -8817

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2547.5
$ws.Range("J16").Value = 1397
$ws.Range("L16").Value = 1397
$ws.Range("N16").Value = -1971

$ws.Range("H22").Value = 1177.1923
$ws.Range("I22").Value = 1011.5
$ws.Range("K22").Value = 1011.5
$ws.Range("M22").Value = -661.5

$ws.Range("H62").Value = 12256.583
$ws.Range("I62").Value = 12530.333
$ws.Range("J62").Value = 11982.833
$ws.Range("K62").Value = 12530.333
$ws.Range("L62").Value = 11982.833
$ws.Range("M62").Value = -11906.333
$ws.Range("N62").Value = -13230.833

$ws.Range("H65").Value = 12256.583
$ws.Range("I65").Value = 12530.333
$ws.Range("J65").Value = 11982.833
$ws.Range("K65").Value = 62651.665
$ws.Range("L65").Value = 59914.165
$ws.Range("M65").Value = -59531.665
$ws.Range("N65").Value = -66154.16500000001

$ws.Range("H113").Value = 2547.5
$ws.Range("J113").Value = 1397
$ws.Range("L113").Value = 1397
$ws.Range("N113").Value = -5737

$ws.Range("H134").Value = 2409947.8
$ws.Range("I134").Value = 2610535
$ws.Range("J134").Value = 2898.5
$ws.Range("K134").Value = 7831605
$ws.Range("L134").Value = 8695.5
$ws.Range("M134").Value = -7829070
$ws.Range("N134").Value = -13765.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 25324.4
$ws.Range("J97").Value = 1236.8
$ws.Range("L97").Value = 3710.4
$ws.Range("N97").Value = -4702.4

$ws.Range("H132").Value = 9288413
$ws.Range("I132").Value = 1099.8572
$ws.Range("J132").Value = 15198522
$ws.Range("K132").Value = 9898.7148
$ws.Range("L132").Value = 136786698
$ws.Range("M132").Value = -7368.7148
$ws.Range("N132").Value = -136791758

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 11206.689
$ws.Range("I102").Value = 14694.95
$ws.Range("J102").Value = 3455
$ws.Range("K102").Value = 14694.95
$ws.Range("L102").Value = 3455
$ws.Range("M102").Value = -13072.95
$ws.Range("N102").Value = -6699

$ws.Range("H113").Value = 4248.1665
$ws.Range("J113").Value = 3497.8
$ws.Range("L113").Value = 3497.8
$ws.Range("N113").Value = -7837.8

$ws.Range("H132").Value = 2923.5854
$ws.Range("I132").Value = 2253.1282
$ws.Range("K132").Value = 6759.3846
$ws.Range("M132").Value = -4229.3846

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 39752.637
$ws.Range("I40").Value = 47860.5
$ws.Range("K40").Value = 47860.5
$ws.Range("M40").Value = -47724.5

$ws.Range("H46").Value = 2201.2778
$ws.Range("J46").Value = 2599.6
$ws.Range("L46").Value = 2599.6
$ws.Range("N46").Value = -2975.6

$ws.Range("H61").Value = 1693.5555
$ws.Range("J61").Value = 3812.25
$ws.Range("L61").Value = 3812.25
$ws.Range("N61").Value = -4216.25

$ws.Range("H68").Value = 2812.5652
$ws.Range("I68").Value = 1638.2858
$ws.Range("J68").Value = 4639.222
$ws.Range("K68").Value = 1638.2858
$ws.Range("L68").Value = 4639.222
$ws.Range("M68").Value = -889.2858000000001
$ws.Range("N68").Value = -6137.222

$ws.Range("H71").Value = 2812.5652
$ws.Range("I71").Value = 1638.2858
$ws.Range("J71").Value = 4639.222
$ws.Range("K71").Value = 8191.429
$ws.Range("L71").Value = 23196.11
$ws.Range("M71").Value = -4447.429
$ws.Range("N71").Value = -30684.11

$ws.Range("H100").Value = 6031.8237
$ws.Range("J100").Value = 3291.3333
$ws.Range("L100").Value = 3291.3333
$ws.Range("N100").Value = -4373.3333

$ws.Range("H113").Value = 1693.5555
$ws.Range("J113").Value = 3812.25
$ws.Range("L113").Value = 3812.25
$ws.Range("N113").Value = -8152.25

$ws.Range("H136").Value = 4600
$ws.Range("I136").Value = 2081.5454
$ws.Range("K136").Value = 6244.6362
$ws.Range("M136").Value = -3694.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 95057
$ws.Range("I64").Value = 90000
$ws.Range("K64").Value = 90000
$ws.Range("M64").Value = -89752

$ws.Range("H67").Value = 95057
$ws.Range("I67").Value = 90000
$ws.Range("K67").Value = 90000
$ws.Range("M67").Value = -89142

$ws.Range("H113").Value = 2634.4614
$ws.Range("I113").Value = 876.8823
$ws.Range("J113").Value = 5954.3335
$ws.Range("K113").Value = 2630.6469
$ws.Range("L113").Value = 17863.0005
$ws.Range("M113").Value = -460.6468999999997
$ws.Range("N113").Value = -22203.0005

$ws.Range("H122").Value = 6185.6665
$ws.Range("I122").Value = 3945.4443
$ws.Range("K122").Value = 11836.3329
$ws.Range("M122").Value = -9386.332900000001

$ws.Range("H136").Value = 2786.3953
$ws.Range("I136").Value = 2016.1842
$ws.Range("J136").Value = 8640
$ws.Range("K136").Value = 6048.5526
$ws.Range("L136").Value = 25920
$ws.Range("M136").Value = -3498.5526
$ws.Range("N136").Value = -31020

